$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "max" column (C) is removed entirely; the following columns
# ("prediction" and "rejection-f") shift one position to the left,
# becoming the new C and D columns, and the sheet's used range shrinks
# from A1:E2 to A1:D2.
$ws.Columns.Item(3).Delete()

# The single remaining child's predicted confidence value is updated.
$ws.Range("B2").Value = 98.93711889947755
